# Apply cryptocurrency price/volume updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain/pure numeric string (e.g. "569.03") must be
# forced to Text format first, otherwise Excel auto-converts them into numbers -
# the source data stores these as literal text (inline strings), e.g. "24.09" not 24.09.
$textCells = @("D5", "D12", "D13", "D17", "D20", "D21", "D25", "D28", "D31", "D32", "D33", "D34", "D36", "D39", "D44", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# New values for the "Price" (D) and "Volume(1h)" (E) columns
$updates = @{
    'D5' = '569.03'
    'D12' = '0.353'
    'D13' = '4.87'
    'D17' = '24.09'
    'D20' = '7.37'
    'D21' = '346.61'
    'D25' = '69.26'
    'D28' = '8.60'
    'D31' = '7.54'
    'D32' = '437.78'
    'D33' = '1.18'
    'D34' = '0.999'
    'D36' = '156.45'
    'D39' = '18.11'
    'D44' = '2.36'
    'D50' = '0.0722'
    'D51' = '0.571'
    'D2' = '69.471.57'
    'E2' = '  +0.26%  '
    'D3' = '2.489.29'
    'E3' = '  -0.93%  '
    'E4' = '  -0.02%  '
    'E5' = '  -0.57%  '
    'E6' = '  -1.51%  '
    'E7' = '  -0.05%  '
    'E8' = '  -0.86%  '
    'E9' = '  -0.89%  '
    'E10' = '  -1.41%  '
    'E11' = '  -0.57%  '
    'E12' = '  -1.23%  '
    'E13' = '  -1.02%  '
    'D14' = '2.944.71'
    'D15' = '69.330.29'
    'E15' = '  -0.01%  '
    'E16' = '  -0.52%  '
    'E17' = '  -3.14%  '
    'D18' = '2.491.87'
    'E18' = '  -0.97%  '
    'E19' = '  -1.83%  '
    'E20' = '  -4.28%  '
    'E21' = '  -0.76%  '
    'E22' = '  -1.51%  '
    'E23' = '  -4.16%  '
    'E24' = '  +0.02%  '
    'E25' = '  -1.10%  '
    'E26' = '  -2.89%  '
    'D27' = '2.616.32'
    'E27' = '  -1.26%  '
    'E28' = '  -3.62%  '
    'E29' = '  +0.68%  '
    'E30' = '  -3.14%  '
    'E31' = '  -4.25%  '
    'E33' = '  -4.58%  '
    'E34' = '  +0.03%  '
    'E36' = '  -0.76%  '
    'E37' = '  -3.23%  '
    'E38' = '  +0.29%  '
    'E39' = '  -2.35%  '
    'E41' = '  -2.16%  '
    'E42' = '  -4.12%  '
    'E43' = '  -2.17%  '
    'E44' = '  +47.56%  '
    'E45' = '  -5.11%  '
    'E46' = '  -6.08%  '
    'E47' = '  -2.73%  '
    'E48' = '  -2.12%  '
    'E49' = '  -4.45%  '
    'E50' = '  -0.99%  '
    'E51' = '  -0.88%  '
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
